$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "37.855.51"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.083.71"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.30"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.41%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.396"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").Value = "2.070.84"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "37.776.00"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  +4.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("E26").Value = "  +3.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.136"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  +2.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.16%  "
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.41%  "
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.46%  "
$ws.Range("D44").Value = "1.446.72"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.50%  "
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "2.275.85"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.34%  "
